$wb = $excel.ActiveWorkbook

# Add a new worksheet "AWS" after the last existing sheet ("Django")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "AWS"

# Match column widths used by the other course sheets (title / url columns)
$ws.Columns.Item(2).ColumnWidth = 82.28515625
$ws.Columns.Item(3).ColumnWidth = 118.85546875

# Row 2 (write URL before Title so shared-string indices come out URL, Title)
$ws.Range("C2").Value = "https://www.udemy.com/course/networking-in-aws/"
$ws.Range("B2").Value = "AWS VPC and Networking in depth: Learn practically in 8 hrs"

# Row 4
$ws.Range("C4").Value = "https://www.udemy.com/course/aws-with-python-and-boto3-managing-ec2-and-vpc/"
$ws.Range("B4").Value = "Managing EC2 and VPC: AWS with Python and Boto3 Series"

# Row 6
$ws.Range("C6").Value = "https://www.udemy.com/course/aws-automation-with-boto3-of-python-and-lambda-functions/"
$ws.Range("B6").Value = "AWS Automation with boto3 of Python and Lambda Functions"

# Row 8
$ws.Range("C8").Value = "https://www.udemy.com/course/aws-ec2-masterclass/"
$ws.Range("B8").Value = "Amazon EC2 Master Class (with Auto Scaling & Load Balancer)"

# Row 10
$ws.Range("C10").Value = "https://www.udemy.com/course/hosting-websites-with-amazon-lightsail/"
$ws.Range("B10").Value = "Hosting Websites with Amazon Lightsail | AWS for 2020"

# Match the selection state used by the other sheets (last populated cell selected)
$ws.Range("B10").Select() | Out-Null
